$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 38 / 39 — the two match records were re-ordered (their full data,
#    columns B through AC, swap places; column A "id" stays put per row).
# ---------------------------------------------------------------------------
$ws.Range("B38").Value = 6782522
$ws.Range("F38").Value = "Puntarenas"
$ws.Range("G38").Value = "Sporting San Jose"
$ws.Range("H38").Value = 1
$ws.Range("I38").Value = 2
$ws.Range("J38").Value = "A"
$ws.Range("K38").Value = 2.5
$ws.Range("L38").Value = 3.5
$ws.Range("M38").Value = 2.5
$ws.Range("N38").Value = 2.2
$ws.Range("O38").Value = 3.5
$ws.Range("P38").Value = 2.9
$ws.Range("Q38").Value = -0.25
$ws.Range("R38").Value = 1.9
$ws.Range("S38").Value = 1.9
$ws.Range("T38").Value = 2.5
$ws.Range("U38").Value = 1.9
$ws.Range("V38").Value = 1.9
$ws.Range("W38").Value = -1
$ws.Range("X38").Value = -1
$ws.Range("Y38").Value = 1.9
$ws.Range("Z38").Value = -1
$ws.Range("AA38").Value = 0.8999999999999999
$ws.Range("AB38").Value = 0.8999999999999999
$ws.Range("AC38").Value = -1

$ws.Range("B39").Value = 6781354
$ws.Range("F39").Value = "Puntarenas"
$ws.Range("G39").Value = "AD San Carlos"
$ws.Range("H39").Value = 1
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = "D"
$ws.Range("K39").Value = 2.4
$ws.Range("L39").Value = 3.2
$ws.Range("M39").Value = 2.8
$ws.Range("N39").Value = 2.3
$ws.Range("O39").Value = 3.2
$ws.Range("P39").Value = 3
$ws.Range("Q39").Value = -0.25
$ws.Range("R39").Value = 2
$ws.Range("S39").Value = 1.8
$ws.Range("T39").Value = 2.25
$ws.Range("U39").Value = 1.9
$ws.Range("V39").Value = 1.9
$ws.Range("W39").Value = 1.3
$ws.Range("X39").Value = -1
$ws.Range("Y39").Value = -1
$ws.Range("Z39").Value = 1
$ws.Range("AA39").Value = -1
$ws.Range("AB39").Value = -1
$ws.Range("AC39").Value = 0.8999999999999999

# ---------------------------------------------------------------------------
# 2) Rows 110 / 111 — same kind of swap.
# ---------------------------------------------------------------------------
$ws.Range("B110").Value = 6782581
$ws.Range("F110").Value = "Alajuelense"
$ws.Range("G110").Value = "AD Grecia"
$ws.Range("H110").Value = 2
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = "D"
$ws.Range("K110").Value = 1.181
$ws.Range("L110").Value = 6.5
$ws.Range("M110").Value = 11
$ws.Range("N110").Value = 1.25
$ws.Range("O110").Value = 5
$ws.Range("P110").Value = 9
$ws.Range("Q110").Value = -1.75
$ws.Range("R110").Value = 1.975
$ws.Range("S110").Value = 1.825
$ws.Range("T110").Value = 3.25
$ws.Range("U110").Value = 2
$ws.Range("V110").Value = 1.8
$ws.Range("W110").Value = 0.25
$ws.Range("X110").Value = -1
$ws.Range("Y110").Value = -1
$ws.Range("Z110").Value = 0.4875
$ws.Range("AA110").Value = -0.5
$ws.Range("AB110").Value = -1
$ws.Range("AC110").Value = 0.8

$ws.Range("B111").Value = 6782579
$ws.Range("F111").Value = "Municipal Perez Zeledon"
$ws.Range("G111").Value = "AD San Carlos"
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 2
$ws.Range("J111").Value = "A"
$ws.Range("K111").Value = 2.4
$ws.Range("L111").Value = 3.3
$ws.Range("M111").Value = 2.7
$ws.Range("N111").Value = 2.375
$ws.Range("O111").Value = 3.4
$ws.Range("P111").Value = 2.8
$ws.Range("Q111").Value = -0.25
$ws.Range("R111").Value = 2
$ws.Range("S111").Value = 1.8
$ws.Range("T111").Value = 2.5
$ws.Range("U111").Value = 1.875
$ws.Range("V111").Value = 1.925
$ws.Range("W111").Value = -1
$ws.Range("X111").Value = -1
$ws.Range("Y111").Value = 1.8
$ws.Range("Z111").Value = -1
$ws.Range("AA111").Value = 0.8
$ws.Range("AB111").Value = -1
$ws.Range("AC111").Value = 0.925

# ---------------------------------------------------------------------------
# 3) Rows 130 / 131 — same kind of swap.
# ---------------------------------------------------------------------------
$ws.Range("B130").Value = 6782596
$ws.Range("F130").Value = "Herediano"
$ws.Range("G130").Value = "AD Guanacasteca"
$ws.Range("H130").Value = 3
$ws.Range("I130").Value = 4
$ws.Range("J130").Value = "A"
$ws.Range("K130").Value = 1.363
$ws.Range("L130").Value = 4.75
$ws.Range("M130").Value = 8
$ws.Range("N130").Value = 1.444
$ws.Range("O130").Value = 4.333
$ws.Range("P130").Value = 7
$ws.Range("Q130").Value = -1.25
$ws.Range("R130").Value = 1.975
$ws.Range("S130").Value = 1.825
$ws.Range("T130").Value = 2.75
$ws.Range("U130").Value = 1.775
$ws.Range("V130").Value = 2.025
$ws.Range("W130").Value = -1
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = 6
$ws.Range("Z130").Value = -1
$ws.Range("AA130").Value = 0.825
$ws.Range("AB130").Value = 0.7749999999999999
$ws.Range("AC130").Value = -1

$ws.Range("B131").Value = 6782595
$ws.Range("F131").Value = "Santos de Gupiles"
$ws.Range("G131").Value = "Sporting San Jose"
$ws.Range("H131").Value = 3
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = "D"
$ws.Range("K131").Value = 1.4
$ws.Range("L131").Value = 4.75
$ws.Range("M131").Value = 7
$ws.Range("N131").Value = 1.363
$ws.Range("O131").Value = 4.75
$ws.Range("P131").Value = 8.5
$ws.Range("Q131").Value = -1.25
$ws.Range("R131").Value = 1.8
$ws.Range("S131").Value = 2
$ws.Range("T131").Value = 3
$ws.Range("U131").Value = 1.95
$ws.Range("V131").Value = 1.85
$ws.Range("W131").Value = 0.363
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = 0.8
$ws.Range("AA131").Value = -1
$ws.Range("AB131").Value = 0
$ws.Range("AC131").Value = -0

# ---------------------------------------------------------------------------
# 4) New rows 228-231 — four newly scraped fixtures appended at the bottom.
#    Copy the number-format / border styling from the last existing row (227)
#    onto column A (id style) and column E (date style) of each new row so
#    the same style indices (s="1" / s="2") get reused.
# ---------------------------------------------------------------------------
$ws.Range("A227").Copy()
$ws.Range("A228").PasteSpecial(-4122)
$ws.Range("A229").PasteSpecial(-4122)
$ws.Range("A230").PasteSpecial(-4122)
$ws.Range("A231").PasteSpecial(-4122)

$ws.Range("E227").Copy()
$ws.Range("E228").PasteSpecial(-4122)
$ws.Range("E229").PasteSpecial(-4122)
$ws.Range("E230").PasteSpecial(-4122)
$ws.Range("E231").PasteSpecial(-4122)

# Row 228
$ws.Range("A228").Value = 226
$ws.Range("B228").Value = 8048492
$ws.Range("C228").Value = "Costa Rica Primera Division"
$ws.Range("D228").Value = "Costa Rica Primera Division"
$ws.Range("E228").Value = 45389.75
$ws.Range("F228").Value = "AD Guanacasteca"
$ws.Range("G228").Value = "Puntarenas"
$ws.Range("H228").Value = 3
$ws.Range("I228").Value = 0
$ws.Range("J228").Value = "H"
$ws.Range("K228").Value = 2.2
$ws.Range("L228").Value = 3.2
$ws.Range("M228").Value = 3.3
$ws.Range("N228").Value = 2.2
$ws.Range("O228").Value = 3.2
$ws.Range("P228").Value = 3.3
$ws.Range("Q228").Value = -0.25
$ws.Range("R228").Value = 1.875
$ws.Range("S228").Value = 1.925
$ws.Range("T228").Value = 2.25
$ws.Range("U228").Value = 2
$ws.Range("V228").Value = 1.8
$ws.Range("W228").Value = 1.2
$ws.Range("X228").Value = -1
$ws.Range("Y228").Value = -1
$ws.Range("Z228").Value = 0.875
$ws.Range("AA228").Value = -1
$ws.Range("AB228").Value = 1
$ws.Range("AC228").Value = -1

# Row 229
$ws.Range("A229").Value = 227
$ws.Range("B229").Value = 7623997
$ws.Range("C229").Value = "Costa Rica Primera Division"
$ws.Range("D229").Value = "Costa Rica Primera Division"
$ws.Range("E229").Value = 45389.83333333334
$ws.Range("F229").Value = "Deportivo Saprissa"
$ws.Range("G229").Value = "Sporting San Jose"
$ws.Range("H229").Value = 2
$ws.Range("I229").Value = 1
$ws.Range("J229").Value = "H"
$ws.Range("K229").Value = 1.285
$ws.Range("L229").Value = 5
$ws.Range("M229").Value = 8
$ws.Range("N229").Value = 1.285
$ws.Range("O229").Value = 5.25
$ws.Range("P229").Value = 7.5
$ws.Range("Q229").Value = -1.5
$ws.Range("R229").Value = 1.9
$ws.Range("S229").Value = 1.9
$ws.Range("T229").Value = 2.75
$ws.Range("U229").Value = 1.8
$ws.Range("V229").Value = 2
$ws.Range("W229").Value = 0.2849999999999999
$ws.Range("X229").Value = -1
$ws.Range("Y229").Value = -1
$ws.Range("Z229").Value = -1
$ws.Range("AA229").Value = 0.8999999999999999
$ws.Range("AB229").Value = 0.4
$ws.Range("AC229").Value = -0.5

# Row 230 (match not played yet - no FTHG/FTAG/FTR, no PL_AhOver/PL_AhUnder)
$ws.Range("A230").Value = 228
$ws.Range("B230").Value = 8070744
$ws.Range("C230").Value = "Costa Rica Primera Division"
$ws.Range("D230").Value = "Costa Rica Primera Division"
$ws.Range("E230").Value = 45394.95833333334
$ws.Range("F230").Value = "Municipal Liberia"
$ws.Range("G230").Value = "AD Guanacasteca"
$ws.Range("K230").Value = 1.833
$ws.Range("L230").Value = 3.5
$ws.Range("M230").Value = 4
$ws.Range("N230").Value = 1.833
$ws.Range("O230").Value = 3.5
$ws.Range("P230").Value = 4
$ws.Range("Q230").Value = -0.5
$ws.Range("R230").Value = 1.825
$ws.Range("S230").Value = 1.975
$ws.Range("T230").Value = 2.5
$ws.Range("U230").Value = 1.9
$ws.Range("V230").Value = 1.9
$ws.Range("W230").Value = 0
$ws.Range("X230").Value = 0
$ws.Range("Y230").Value = 0
$ws.Range("Z230").Value = 0
$ws.Range("AA230").Value = 0

# Row 231 (match not played yet - no FTHG/FTAG/FTR, no PL_AhOver/PL_AhUnder)
$ws.Range("A231").Value = 229
$ws.Range("B231").Value = 7623999
$ws.Range("C231").Value = "Costa Rica Primera Division"
$ws.Range("D231").Value = "Costa Rica Primera Division"
$ws.Range("E231").Value = 45395.79166666666
$ws.Range("F231").Value = "Puntarenas"
$ws.Range("G231").Value = "Sporting San Jose"
$ws.Range("K231").Value = 2.15
$ws.Range("L231").Value = 3.4
$ws.Range("M231").Value = 3.1
$ws.Range("N231").Value = 2.15
$ws.Range("O231").Value = 3.4
$ws.Range("P231").Value = 3.1
$ws.Range("Q231").Value = -0.25
$ws.Range("R231").Value = 1.875
$ws.Range("S231").Value = 1.925
$ws.Range("T231").Value = 2.5
$ws.Range("U231").Value = 2.025
$ws.Range("V231").Value = 1.775
$ws.Range("W231").Value = 0
$ws.Range("X231").Value = 0
$ws.Range("Y231").Value = 0
$ws.Range("Z231").Value = 0
$ws.Range("AA231").Value = 0
